$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Header row (B1:G1): "MicroclustersNumberNN" -> "micro-cluster=NN" ---
$ws.Range("B1").Value = "micro-cluster=25"
$ws.Range("C1").Value = "micro-cluster=50"
$ws.Range("D1").Value = "micro-cluster=75"
$ws.Range("E1").Value = "micro-cluster=100"
$ws.Range("F1").Value = "micro-cluster=125"
$ws.Range("G1").Value = "micro-cluster=150"

# --- Data rows 2-7 (B:G): normalize throughput by dividing by 1000 ---
$ws.Range("B2").Value = 4368.0540149999997 / 1000
$ws.Range("C2").Value = 3517.8845379999998 / 1000
$ws.Range("D2").Value = 3018.7886189999999 / 1000
$ws.Range("E2").Value = 2498.587039 / 1000
$ws.Range("F2").Value = 2304.1981460000002 / 1000
$ws.Range("G2").Value = 1930.053005 / 1000

$ws.Range("B3").Value = 8136.9243669999996 / 1000
$ws.Range("C3").Value = 6669.5220010000003 / 1000
$ws.Range("D3").Value = 5728.6454199999998 / 1000
$ws.Range("E3").Value = 4971.6586530000004 / 1000
$ws.Range("F3").Value = 4472.9249099999997 / 1000
$ws.Range("G3").Value = 3981.4481839999999 / 1000

$ws.Range("B4").Value = 14847.83152 / 1000
$ws.Range("C4").Value = 12735.207689999999 / 1000
$ws.Range("D4").Value = 10890.130020000001 / 1000
$ws.Range("E4").Value = 9410.8588099999997 / 1000
$ws.Range("F4").Value = 8514.0148819999995 / 1000
$ws.Range("G4").Value = 7639.1553599999997 / 1000

$ws.Range("B5").Value = 26673.65033 / 1000
$ws.Range("C5").Value = 23196.054550000001 / 1000
$ws.Range("D5").Value = 20574.95535 / 1000
$ws.Range("E5").Value = 14096.87846 / 1000
$ws.Range("F5").Value = 12413.46162 / 1000
$ws.Range("G5").Value = 10421.524359999999 / 1000

$ws.Range("B6").Value = 40895.660179999999 / 1000
$ws.Range("C6").Value = 36844.791279999998 / 1000
$ws.Range("D6").Value = 30373.753546248699 / 1000
$ws.Range("E6").Value = 24477.873553169899 / 1000
$ws.Range("F6").Value = 21331.428400000001 / 1000
$ws.Range("G6").Value = 19046.62617 / 1000

$ws.Range("B7").Value = 48094.043080000003 / 1000
$ws.Range("C7").Value = 46481.566899999998 / 1000
$ws.Range("D7").Value = 40110.365010000001 / 1000
$ws.Range("E7").Value = 33539.271468442603 / 1000
$ws.Range("F7").Value = 29585.540468250201 / 1000
$ws.Range("G7").Value = 26435.5785935493 / 1000

# --- Update the selected cell to match the new view state ---
$ws.Range("C12").Select()
